$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to reflect new TPM-derived specificity scores
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Remove row 3 (the MuSCs target-cluster record no longer present after rerun)
$ws.Rows.Item(3).Delete()
